$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.939.49"
$ws.Range("E2").Value = "  +0.83%  "

$ws.Range("D3").Value = "2.631.06"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'596.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").Value = "'153.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.43%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -0.96%  "

$ws.Range("D9").Value = "2.631.14"
$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("E10").Value = "  +9.27%  "

$ws.Range("E11").Value = "  -0.70%  "

$ws.Range("E12").Value = "  +0.44%  "

$ws.Range("D13").Value = "'0.348"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.06%  "

$ws.Range("D14").Value = "'27.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("E15").Value = "  +3.57%  "

$ws.Range("D16").Value = "3.111.38"
$ws.Range("E16").Value = "  +0.57%  "

$ws.Range("D17").Value = "67.808.33"
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("D18").Value = "2.623.26"
$ws.Range("E18").Value = "  -0.14%  "

$ws.Range("D19").Value = "'11.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.71%  "

$ws.Range("D20").Value = "'371.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.28%  "

$ws.Range("E21").Value = "  +0.10%  "

$ws.Range("E22").Value = "  -1.05%  "

$ws.Range("E23").Value = "  -1.88%  "

$ws.Range("E24").Value = "  -0.74%  "

$ws.Range("D25").Value = "'72.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.28%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").Value = "'9.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.99%  "

$ws.Range("D28").Value = "2.766.29"

$ws.Range("E29").Value = "  +2.06%  "

$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("D31").Value = "'577.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.90%  "

$ws.Range("E32").Value = "  -0.48%  "

$ws.Range("D33").Value = "'7.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.55%  "

$ws.Range("E34").Value = "  +0.54%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").Value = "'0.127"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.35%  "

$ws.Range("D37").Value = "'1.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("D38").Value = "'158.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.67%  "

$ws.Range("D39").Value = "'19.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.37%  "

$ws.Range("E40").Value = "  +4.82%  "

$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("E42").Value = "  +1.61%  "

$ws.Range("D43").Value = "0.0₆0338"
$ws.Range("E43").Value = "  +16.44%  "

$ws.Range("D44").Value = "'2.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.81%  "

$ws.Range("D45").Value = "'17.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.24%  "

$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").Value = "'40.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.19%  "

$ws.Range("D48").Value = "'156.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.55%  "

$ws.Range("D49").Value = "'3.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.88%  "

$ws.Range("D50").Value = "'21.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.66%  "

$ws.Range("E51").Value = "  -1.30%  "
